# Update "想去人数" (number of people wanting to go) values for the
# three upcoming events listed on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 84
    $ws.Range("F3").Value = 334
    $ws.Range("F4").Value = 4583
}
